# Append a new "terminal command" style paragraph (dark-shaded, Consolas
# font) after the existing riscv64-unknown-elf-gdb.exe / wsl$ command line,
# plus a trailing blank paragraph — mirrors pasting a second gdb command
# block (Windows/WSL environment) at the very end of the document body.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- New formatted paragraph with the WSL-path gdb command split into runs ---
$paraXml = @"
<w:p $wNs>
  <w:pPr>
    <w:shd w:val="clear" w:color="auto" w:fill="2A2A2A"/>
    <w:spacing w:line="450" w:lineRule="atLeast"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
      <w:color w:val="E2E2E2"/>
      <w:sz w:val="33"/>
      <w:szCs w:val="33"/>
      <w:lang w:eastAsia="en-GB"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
      <w:color w:val="656565"/>
      <w:sz w:val="33"/>
      <w:szCs w:val="33"/>
      <w:lang w:eastAsia="en-GB"/>
    </w:rPr>
    <w:t>riscv64-unknown-elf-gdb.exe Z:/home/imandadras/diana-riscv-src/ana_char_loop/build/hwme.c/pulpissimo/hwme/hwme</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
      <w:color w:val="656565"/>
      <w:sz w:val="33"/>
      <w:szCs w:val="33"/>
      <w:lang w:eastAsia="en-GB"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
      <w:color w:val="656565"/>
      <w:sz w:val="33"/>
      <w:szCs w:val="33"/>
      <w:lang w:eastAsia="en-GB"/>
    </w:rPr>
    <w:t>-x</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
      <w:color w:val="656565"/>
      <w:sz w:val="33"/>
      <w:szCs w:val="33"/>
      <w:lang w:eastAsia="en-GB"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/>
      <w:color w:val="656565"/>
      <w:sz w:val="33"/>
      <w:szCs w:val="33"/>
      <w:lang w:eastAsia="en-GB"/>
    </w:rPr>
    <w:t>C:/zedboard/diana-fpga-sw/host_scripts/templates/gdb-run-soc.sh</w:t>
  </w:r>
</w:p>
"@

$tail = $d.Content
$tail.Collapse(0)            # wdCollapseEnd - position right at end of body, before sectPr
$null = $tail.InsertXML($paraXml)

# --- Trailing blank paragraph after the new command block ---
$tail2 = $d.Content
$tail2.Collapse(0)
$null = $tail2.InsertXML("<w:p $wNs/>")
